$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$abstractText = "Given the ongoing SARS-CoV-2 pandemic, identification of immunogenic targets against the coronavirus spike glycoprotein will provide crucial advances towards the development of sensitive diagnostic tools and potential vaccine candidate targets.`n In this study, using pools of overlapping linear B-cell peptides, we report two IgG immunodominant regions on SARS-CoV-2 spike glycoprotein that are recognised by sera from COVID-19 convalescent patients.`n Notably, one is specific to SARS-CoV-2, which is located in close proximity to the receptor binding domain.`n The other region, which is localised at the fusion peptide, could potentially function as a pan-SARS target.`n Functionally, antibody depletion assays demonstrate that antibodies targeting these immunodominant regions significantly alter virus neutralisation capacities.`n Taken together, identification and validation of these neutralising B-cell epitopes will provide insights towards the design of diagnostics and vaccine candidates against this high priority coronavirus.`n"

$authorsText = "[Chek Meng%Poh%NULL%1,   Guillaume%Carissimo%NULL%1,   Bei%Wang%NULL%2,   Bei%Wang%NULL%0,   Siti Naqiah%Amrun%NULL%1,   Cheryl Yi-Pin%Lee%NULL%1,   Rhonda Sin-Ling%Chee%NULL%1,   Siew-Wai%Fong%NULL%1,   Nicholas Kim-Wah%Yeo%NULL%1,   Wen-Hsin%Lee%NULL%1,   Anthony%Torres-Ruesta%NULL%1,   Yee-Sin%Leo%NULL%2,   Yee-Sin%Leo%NULL%0,   Mark I-Cheng%Chen%NULL%1,   Seow-Yen%Tan%NULL%1,   Louis Yi Ann%Chai%NULL%1,   Shirin%Kalimuddin%NULL%1,   Shirley Seah Gek%Kheng%NULL%1,   Siew-Yee%Thien%NULL%1,   Barnaby Edward%Young%NULL%1,   David C.%Lye%NULL%1,   Brendon John%Hanson%NULL%1,   Cheng-I%Wang%NULL%1,   Laurent%Renia%NULL%1,   Lisa F. P.%Ng%lisa_ng@immunol.a-star.edu.sg%2,   Lisa F. P.%Ng%lisa_ng@immunol.a-star.edu.sg%0]"

$ws.Range("D7").Value = $abstractText
$ws.Range("E7").Value = $authorsText

$ws.Rows.Item(7).EntireRow.AutoFit()
